# Updated cryptos list on Sat Jul  6 11:56:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.726.48'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.007.20'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.03%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '513.48'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +5.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.84'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +6.85%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.435'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +4.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.53'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +7.83%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +9.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.357'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.54%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.520.38'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.77'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +5.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000159'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +14.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '56.761.17'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.003.10'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.98'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +6.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.55'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.86'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +6.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '331.59'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +6.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.484'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.19'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +5.65%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +10.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0₃0910'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +7.99%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.75'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.05'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +9.28%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +8.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.82'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +8.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.75'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +8.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '153.85'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.58'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.69'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.15%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.87%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.77'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.037.17'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.90'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.22%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.281.05'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +8.16%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.69'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +5.13%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.24%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +14.30%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.85'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.47'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0873'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +6.47%  '
